$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (Paris Saint-Germain match), shifting
# the existing rows 4-29 down by one.
$ws.Rows.Item(4).Insert()

$ws.Cells.Item(4, 1).Value = "Fri Oct 17"
$ws.Cells.Item(4, 2).Value = "Paris Saint-Germain  - RC Strasbourg Alsace: 3:3"
$ws.Cells.Item(4, 3).Value = 2.72
$ws.Cells.Item(4, 4).Value = "Paris Saint-Germain"
$ws.Cells.Item(4, 5).Value = 3.5
$ws.Cells.Item(4, 6).Value = "'73%"
$ws.Cells.Item(4, 8).Value = 6
$ws.Cells.Item(4, 9).Value = $false

# Insert a second new row at position 26 (Llandudno FC match), shifting
# the rows that follow (Sanfrecce Hiroshima, MFK Dukla Banska Bystrica and
# the summary formulas) down by one more.
$ws.Rows.Item(26).Insert()

$ws.Cells.Item(26, 1).Value = "Fri Oct 17"
$ws.Cells.Item(26, 2).Value = "Llandudno FC ✓ - Flint Mountain: 3:0"
$ws.Cells.Item(26, 3).Value = 3.24
$ws.Cells.Item(26, 4).Value = "Llandudno FC"
$ws.Cells.Item(26, 5).Value = 4.5
$ws.Cells.Item(26, 6).Value = "'55%"
$ws.Cells.Item(26, 7).Value = "✓"
$ws.Cells.Item(26, 8).Value = 3
$ws.Cells.Item(26, 9).Value = $true
